# Auto-generated script applying the scheduled-runner price update diff
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC (43 cell updates) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("M6").Value = -738.0000200000001
$ws.Range("J6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("N6").ClearContents()
$ws.Range("H6").Value = 283.33334
$ws.Range("K6").Value = 850.0000200000001
$ws.Range("I6").Value = 283.33334
$ws.Range("H33").Value = 130.3
$ws.Range("I33").Value = 128.11111
$ws.Range("K33").Value = 128.11111
$ws.Range("M33").Value = 100.88889
$ws.Range("K62").Value = 36332.832
$ws.Range("H62").Value = 39332.875
$ws.Range("I62").Value = 36332.832
$ws.Range("M62").Value = -35708.832
$ws.Range("H65").Value = 39332.875
$ws.Range("K65").Value = 181664.16
$ws.Range("I65").Value = 36332.832
$ws.Range("M65").Value = -178544.16
$ws.Range("I69").Value = 0
$ws.Range("M69").ClearContents()
$ws.Range("K69").Value = 0
$ws.Range("H69").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("M72").ClearContents()
$ws.Range("H72").Value = 0
$ws.Range("M125").Value = -6009
$ws.Range("H125").Value = 1043.5
$ws.Range("I125").Value = 941
$ws.Range("K125").Value = 8469
$ws.Range("H132").Value = 301568.44
$ws.Range("I132").Value = 334068.53
$ws.Range("K132").Value = 1002205.59
$ws.Range("M132").Value = -999675.5900000001
$ws.Range("M137").Value = -1991.5386
$ws.Range("H137").Value = 2487.7778
$ws.Range("I137").Value = 1513.8462
$ws.Range("K137").Value = 4541.5386
$ws.Range("M141").Value = 1722.7502
$ws.Range("H141").Value = 1084.2142
$ws.Range("K141").Value = 3457.2498
$ws.Range("I141").Value = 1152.4166

# --- Sheet: ARM (36 cell updates) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("M2").Value = -2726.2632
$ws.Range("J2").Value = 3637.2104
$ws.Range("N2").Value = -3863.2104
$ws.Range("L2").Value = 3637.2104
$ws.Range("H2").Value = 3238.2368
$ws.Range("I2").Value = 2839.2632
$ws.Range("K2").Value = 2839.2632
$ws.Range("J32").Value = 43207.668
$ws.Range("L32").Value = 43207.668
$ws.Range("N32").Value = -43781.668
$ws.Range("H32").Value = 37823.594
$ws.Range("K32").Value = 33785.543
$ws.Range("I32").Value = 33785.543
$ws.Range("M32").Value = -33498.543
$ws.Range("H74").Value = 2553842
$ws.Range("K74").Value = 3291146.5
$ws.Range("I74").Value = 3291146.5
$ws.Range("M74").Value = -3290272.5
$ws.Range("H77").Value = 2553842
$ws.Range("K77").Value = 16455732.5
$ws.Range("I77").Value = 3291146.5
$ws.Range("M77").Value = -16451364.5
$ws.Range("M102").Value = -1179.2273
$ws.Range("J102").Value = 4265
$ws.Range("L102").Value = 4265
$ws.Range("N102").Value = -7509
$ws.Range("K102").Value = 2801.2273
$ws.Range("H102").Value = 3026.423
$ws.Range("I102").Value = 2801.2273
$ws.Range("L116").Value = 3637.2104
$ws.Range("H116").Value = 3238.2368
$ws.Range("I116").Value = 2839.2632
$ws.Range("K116").Value = 2839.2632
$ws.Range("M116").Value = -545.2631999999999
$ws.Range("J116").Value = 3637.2104
$ws.Range("N116").Value = -8225.2104

# --- Sheet: BSM (11 cell updates) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("J3").Value = 3637.2104
$ws.Range("L3").Value = 3637.2104
$ws.Range("N3").Value = -3865.2104
$ws.Range("H3").Value = 3238.2368
$ws.Range("K3").Value = 2839.2632
$ws.Range("I3").Value = 2839.2632
$ws.Range("M3").Value = -2725.2632
$ws.Range("H22").Value = 1232.6666
$ws.Range("K22").Value = 350
$ws.Range("I22").Value = 350
$ws.Range("M22").Value = -177

# --- Sheet: CRP (42 cell updates) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("K17").Value = 70000
$ws.Range("H17").Value = 70000
$ws.Range("I17").Value = 70000
$ws.Range("M17").Value = -69826
$ws.Range("M31").Value = -2680.2104
$ws.Range("J31").Value = 9989.182000000001
$ws.Range("L31").Value = 9989.182000000001
$ws.Range("N31").Value = -10579.182
$ws.Range("K31").Value = 2975.2104
$ws.Range("H31").Value = 6738.8047
$ws.Range("I31").Value = 2975.2104
$ws.Range("I34").Value = 2975.2104
$ws.Range("M34").Value = -2773.2104
$ws.Range("J34").Value = 9989.182000000001
$ws.Range("L34").Value = 9989.182000000001
$ws.Range("N34").Value = -10393.182
$ws.Range("H34").Value = 6738.8047
$ws.Range("K34").Value = 2975.2104
$ws.Range("J41").Value = 62514.5
$ws.Range("L41").Value = 62514.5
$ws.Range("N41").Value = -63370.5
$ws.Range("H41").Value = 62514.5
$ws.Range("J51").Value = 50000
$ws.Range("L51").Value = 50000
$ws.Range("N51").Value = -51472
$ws.Range("H51").Value = 50000
$ws.Range("J60").Value = 129997
$ws.Range("L60").Value = 129997
$ws.Range("N60").Value = -131019
$ws.Range("H60").Value = 84748.5
$ws.Range("L61").Value = 50000
$ws.Range("N61").Value = -50696
$ws.Range("H61").Value = 50000
$ws.Range("J61").Value = 50000
$ws.Range("J95").Value = 93495
$ws.Range("L95").Value = 93495
$ws.Range("N95").Value = -98987
$ws.Range("H95").Value = 93495
$ws.Range("M105").Value = -10935.434
$ws.Range("H105").Value = 11802.091
$ws.Range("K105").Value = 12682.434
$ws.Range("I105").Value = 12682.434

# --- Sheet: CUL (23 cell updates) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("J63").Value = 7197.75
$ws.Range("L63").Value = 21593.25
$ws.Range("N63").Value = -23091.25
$ws.Range("H63").Value = 5323.0835
$ws.Range("J66").Value = 7197.75
$ws.Range("N66").Value = -72267.75
$ws.Range("L66").Value = 64779.75
$ws.Range("H66").Value = 5323.0835
$ws.Range("J93").Value = 5998.6
$ws.Range("L93").Value = 17995.8
$ws.Range("N93").Value = -21739.8
$ws.Range("H93").Value = 5665.5
$ws.Range("M102").Value = -14066
$ws.Range("J102").Value = 10715.143
$ws.Range("L102").Value = 32145.429
$ws.Range("N102").Value = -37013.429
$ws.Range("K102").Value = 16500
$ws.Range("H102").Value = 9556.223
$ws.Range("I102").Value = 5500
$ws.Range("M137").Value = -24897
$ws.Range("H137").Value = 9999
$ws.Range("I137").Value = 9999
$ws.Range("K137").Value = 29997

# --- Sheet: GSM (15 cell updates) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("M2").Value = 36.14286
$ws.Range("J2").Value = 2018.1666
$ws.Range("N2").Value = -2244.1666
$ws.Range("L2").Value = 2018.1666
$ws.Range("H2").Value = 1302.9474
$ws.Range("I2").Value = 76.85714
$ws.Range("K2").Value = 76.85714
$ws.Range("M80").Value = -627752.9
$ws.Range("H80").Value = 458909.53
$ws.Range("I80").Value = 628750.9
$ws.Range("K80").Value = 628750.9
$ws.Range("K83").Value = 3143754.5
$ws.Range("H83").Value = 458909.53
$ws.Range("I83").Value = 628750.9
$ws.Range("M83").Value = -3138762.5

# --- Sheet: LTW (7 cell updates) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("L132").Value = 19932
$ws.Range("N132").Value = -24992
$ws.Range("H132").Value = 1021362.56
$ws.Range("I132").Value = 1386661.2
$ws.Range("K132").Value = 4159983.6
$ws.Range("J132").Value = 6644
$ws.Range("M132").Value = -4157453.6

# --- Sheet: WVR (43 cell updates) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("I8").Value = 35101.5
$ws.Range("M8").Value = -34961.5
$ws.Range("H8").Value = 35101.5
$ws.Range("K8").Value = 35101.5
$ws.Range("H81").Value = 2526.6
$ws.Range("K81").Value = 5244.4444
$ws.Range("I81").Value = 2622.2222
$ws.Range("M81").Value = -4183.4444
$ws.Range("J81").Value = 1666
$ws.Range("L81").Value = 3332
$ws.Range("N81").Value = -5454
$ws.Range("N84").Value = -27268
$ws.Range("H84").Value = 2526.6
$ws.Range("K84").Value = 26222.222
$ws.Range("I84").Value = 2622.2222
$ws.Range("M84").Value = -20918.222
$ws.Range("J84").Value = 1666
$ws.Range("L84").Value = 16660
$ws.Range("K107").Value = 4329.4998
$ws.Range("H107").Value = 1616.92
$ws.Range("I107").Value = 1443.1666
$ws.Range("M107").Value = -2409.4998
$ws.Range("J107").Value = 2063.7144
$ws.Range("L107").Value = 6191.1432
$ws.Range("N107").Value = -10031.1432
$ws.Range("K126").Value = 9079.167000000001
$ws.Range("H126").Value = 3408.6667
$ws.Range("I126").Value = 3026.389
$ws.Range("M126").Value = -6609.167000000001
$ws.Range("J126").Value = 5702.3335
$ws.Range("L126").Value = 17107.0005
$ws.Range("N126").Value = -22047.0005
$ws.Range("H132").Value = 5762505
$ws.Range("I132").Value = 5762505
$ws.Range("K132").Value = 17287515
$ws.Range("M132").Value = -17284985
$ws.Range("K136").Value = 6111942.300000001
$ws.Range("I136").Value = 2037314.1
$ws.Range("M136").Value = -6109392.300000001
$ws.Range("J136").Value = 3600.375
$ws.Range("L136").Value = 10801.125
$ws.Range("N136").Value = -15901.125
$ws.Range("H136").Value = 1051271.1
